$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (47 cell ops) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3030.6667
$ws.Range("J32").Value = 3295.3333
$ws.Range("L32").Value = 3295.3333
$ws.Range("N32").Value = -3947.3333
$ws.Range("H40").Value = 1539.64
$ws.Range("I40").Value = 1534.3914
$ws.Range("K40").Value = 1534.3914
$ws.Range("M40").Value = -1359.3914
$ws.Range("H55").Value = 690.25
$ws.Range("I55").Value = 1005.5
$ws.Range("J55").Value = 375
$ws.Range("K55").Value = 1005.5
$ws.Range("L55").Value = 375
$ws.Range("M55").Value = -791.5
$ws.Range("N55").Value = -803
$ws.Range("H103").Value = 1339.6
$ws.Range("J103").Value = 1424.75
$ws.Range("L103").Value = 4274.25
$ws.Range("N103").Value = -5446.25
$ws.Range("H112").Value = 2597.1875
$ws.Range("I112").Value = 1345
$ws.Range("J112").Value = 3014.5833
$ws.Range("K112").Value = 4035
$ws.Range("L112").Value = 9043.749899999999
$ws.Range("M112").Value = -2927
$ws.Range("N112").Value = -11259.7499
$ws.Range("H125").Value = 2425.5833
$ws.Range("I125").Value = 2509.818
$ws.Range("J125").Value = 1499
$ws.Range("K125").Value = 22588.362
$ws.Range("L125").Value = 13491
$ws.Range("M125").Value = -20128.362
$ws.Range("N125").Value = -18411
$ws.Range("H132").Value = 1938.8
$ws.Range("I132").Value = 1821.1666
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 5463.4998
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -2933.4998
$ws.Range("N132").Value = -14052.5
$ws.Range("H138").Value = 3490.8628
$ws.Range("I138").Value = 2265.7334
$ws.Range("J138").Value = 4001.3333
$ws.Range("K138").Value = 6797.2002
$ws.Range("L138").Value = 12003.9999
$ws.Range("M138").Value = -1657.2002
$ws.Range("N138").Value = -22283.9999

# ---- Sheet: ARM (12 cell ops) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3179.7778
$ws.Range("I45").Value = 2815.7334
$ws.Range("K45").Value = 2815.7334
$ws.Range("M45").Value = -2438.7334
$ws.Range("H122").Value = 31834.3
$ws.Range("I122").Value = 31834.3
$ws.Range("K122").Value = 95502.89999999999
$ws.Range("M122").Value = -93052.89999999999
$ws.Range("H132").Value = 1665.8235
$ws.Range("I132").Value = 1551.7142
$ws.Range("K132").Value = 4655.142599999999
$ws.Range("M132").Value = -2125.142599999999

# ---- Sheet: BSM (12 cell ops) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 13680
$ws.Range("J95").Value = 13680
$ws.Range("L95").Value = 13680
$ws.Range("N95").Value = -19172
$ws.Range("H134").Value = 2999.3333
$ws.Range("I134").Value = 2999.3333
$ws.Range("K134").Value = 8997.999899999999
$ws.Range("M134").Value = -6462.999899999999
$ws.Range("H140").Value = 39999
$ws.Range("J140").Value = 39999
$ws.Range("L140").Value = 39999
$ws.Range("N140").Value = -50359

# ---- Sheet: CRP (31 cell ops) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1051.6
$ws.Range("I31").Value = 857
$ws.Range("K31").Value = 857
$ws.Range("M31").Value = -562
$ws.Range("H32").Value = 13340
$ws.Range("I32").Value = 14010
$ws.Range("K32").Value = 14010
$ws.Range("M32").Value = -13694
$ws.Range("H34").Value = 1051.6
$ws.Range("I34").Value = 857
$ws.Range("K34").Value = 857
$ws.Range("M34").Value = -655
$ws.Range("H58").Value = 2565.2354
$ws.Range("I58").Value = 1870.3334
$ws.Range("K58").Value = 1870.3334
$ws.Range("M58").Value = -1667.3334
$ws.Range("H107").Value = 889
$ws.Range("I107").Value = 798.25
$ws.Range("J107").Value = 1131
$ws.Range("K107").Value = 798.25
$ws.Range("L107").Value = 1131
$ws.Range("M107").Value = 1121.75
$ws.Range("N107").Value = -4971
$ws.Range("H132").Value = 6109.5884
$ws.Range("I132").Value = 6109.5884
$ws.Range("K132").Value = 18328.7652
$ws.Range("M132").Value = -15798.7652
$ws.Range("H136").Value = 2565.2354
$ws.Range("I136").Value = 1870.3334
$ws.Range("K136").Value = 5611.0002
$ws.Range("M136").Value = -3061.0002

# ---- Sheet: CUL (70 cell ops) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7622432.5
$ws.Range("I4").Value = 9234701
$ws.Range("K4").Value = 27704103
$ws.Range("M4").Value = -27703991
$ws.Range("H9").Value = 2495
$ws.Range("I9").Value = 2495
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 7485
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -7261
$ws.Range("N9").ClearContents()
$ws.Range("H29").Value = 646.6
$ws.Range("I29").Value = 152
$ws.Range("K29").Value = 456
$ws.Range("M29").Value = -179
$ws.Range("H40").Value = 73.375
$ws.Range("J40").Value = 155
$ws.Range("L40").Value = 620
$ws.Range("N40").Value = -758
$ws.Range("H64").Value = 1797.4
$ws.Range("J64").Value = 1797.4
$ws.Range("L64").Value = 5392.200000000001
$ws.Range("N64").Value = -5932.200000000001
$ws.Range("H67").Value = 1797.4
$ws.Range("J67").Value = 1797.4
$ws.Range("L67").Value = 5392.200000000001
$ws.Range("N67").Value = -7264.200000000001
$ws.Range("H69").Value = 798.5
$ws.Range("J69").Value = 798.5
$ws.Range("L69").Value = 2395.5
$ws.Range("N69").Value = -4017.5
$ws.Range("H70").Value = 17777.111
$ws.Range("I70").Value = 16332
$ws.Range("K70").Value = 48996
$ws.Range("M70").Value = -48681
$ws.Range("H72").Value = 798.5
$ws.Range("J72").Value = 798.5
$ws.Range("L72").Value = 7186.5
$ws.Range("N72").Value = -15298.5
$ws.Range("H73").Value = 17777.111
$ws.Range("I73").Value = 16332
$ws.Range("K73").Value = 48996
$ws.Range("M73").Value = -47904
$ws.Range("H75").Value = 2432.3333
$ws.Range("J75").Value = 2399
$ws.Range("L75").Value = 7197
$ws.Range("N75").Value = -9193
$ws.Range("H78").Value = 2432.3333
$ws.Range("J78").Value = 2399
$ws.Range("L78").Value = 21591
$ws.Range("N78").Value = -31575
$ws.Range("H111").Value = 306.75
$ws.Range("I111").Value = 306.75
$ws.Range("K111").Value = 920.25
$ws.Range("M111").Value = 2146.75
$ws.Range("H120").Value = 2065
$ws.Range("I120").Value = 2065
$ws.Range("K120").Value = 6195
$ws.Range("M120").Value = -1357
$ws.Range("H133").Value = 17498.625
$ws.Range("I133").Value = 4996.3335
$ws.Range("K133").Value = 14989.0005
$ws.Range("M133").Value = -9929.000499999998
$ws.Range("H134").Value = 860.6
$ws.Range("I134").Value = 716.2308
$ws.Range("J134").Value = 1799
$ws.Range("K134").Value = 2148.6924
$ws.Range("L134").Value = 5397
$ws.Range("M134").Value = 2921.3076
$ws.Range("N134").Value = -15537

# ---- Sheet: GSM (11 cell ops) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2596.818
$ws.Range("I102").Value = 2295.7144
$ws.Range("J102").Value = 3123.75
$ws.Range("K102").Value = 2295.7144
$ws.Range("L102").Value = 3123.75
$ws.Range("M102").Value = -673.7143999999998
$ws.Range("N102").Value = -6367.75
$ws.Range("H122").Value = 2638.9
$ws.Range("I122").Value = 2848.25
$ws.Range("K122").Value = 8544.75
$ws.Range("M122").Value = -6094.75

# ---- Sheet: LTW (29 cell ops) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5434.4614
$ws.Range("I7").Value = 3441.8333
$ws.Range("J7").Value = 7142.4287
$ws.Range("K7").Value = 3441.8333
$ws.Range("L7").Value = 7142.4287
$ws.Range("M7").Value = -3329.8333
$ws.Range("N7").Value = -7366.4287
$ws.Range("H32").Value = 2994.5
$ws.Range("I32").Value = 2994.5
$ws.Range("K32").Value = 2994.5
$ws.Range("M32").Value = -2677.5
$ws.Range("H40").Value = 3153.75
$ws.Range("I40").Value = 2114.6667
$ws.Range("J40").Value = 3777.2
$ws.Range("K40").Value = 2114.6667
$ws.Range("L40").Value = 3777.2
$ws.Range("N40").Value = -4049.2
$ws.Range("M40").Value = -1978.6667
$ws.Range("H122").Value = 5493.089
$ws.Range("I122").Value = 4136.64
$ws.Range("K122").Value = 12409.92
$ws.Range("M122").Value = -9959.920000000002
$ws.Range("H126").Value = 5434.4614
$ws.Range("I126").Value = 3441.8333
$ws.Range("J126").Value = 7142.4287
$ws.Range("K126").Value = 10325.4999
$ws.Range("L126").Value = 21427.2861
$ws.Range("M126").Value = -7855.499899999999
$ws.Range("N126").Value = -26367.2861

# ---- Sheet: WVR (23 cell ops) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 15861
$ws.Range("J45").Value = 15861
$ws.Range("L45").Value = 15861
$ws.Range("N45").Value = -16843
$ws.Range("H70").Value = 88000
$ws.Range("I70").Value = 88000
$ws.Range("K70").Value = 88000
$ws.Range("M70").Value = -87685
$ws.Range("H73").Value = 88000
$ws.Range("I73").Value = 88000
$ws.Range("K73").Value = 88000
$ws.Range("M73").Value = -86908
$ws.Range("H122").Value = 4162.9473
$ws.Range("I122").Value = 3051.7144
$ws.Range("K122").Value = 9155.143199999999
$ws.Range("M122").Value = -6705.143199999999
$ws.Range("H136").Value = 3164.348
$ws.Range("I136").Value = 3091.8
$ws.Range("J136").Value = 3648
$ws.Range("K136").Value = 9275.400000000001
$ws.Range("L136").Value = 10944
$ws.Range("M136").Value = -6725.400000000001
$ws.Range("N136").Value = -16044
